$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2: "Content Placeholder 2" - several bullet text tweaks
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$body = $s2.Shapes.Item(2).TextFrame.TextRange

function Replace-SubText($range, [string]$find, [string]$replacement) {
    $full = $range.Text
    $idx = $full.IndexOf($find)
    if ($idx -ge 0) {
        $sub = $range.Characters($idx + 1, $find.Length)
        $sub.Text = $replacement
    }
}

function Strike-SubText($range, [string]$find) {
    $full = $range.Text
    $idx = $full.IndexOf($find)
    if ($idx -ge 0) {
        $sub = $range.Characters($idx + 1, $find.Length)
        $sub.Font.Strikethrough = $true
    }
}

# "Canonical schema (JSON + CDDL)" -> "Canonical schema (JSON + CDDL) Framework Schema"
Replace-SubText $body "Canonical schema (JSON + CDDL)" "Canonical schema (JSON + CDDL) Framework Schema"

# "Canonical Document" gains strikethrough formatting
Strike-SubText $body "Canonical Document"

# "Alt Schema (JSON + CDDL)" -> "Alt Schema (JSON + CDDL) – Need to Merge"
Replace-SubText $body "Alt Schema (JSON + CDDL)" "Alt Schema (JSON + CDDL) – Need to Merge"

# "Alt Document" gains strikethrough formatting
Strike-SubText $body "Alt Document"

# "I-D in progress" -> "I-D in progress – stable (IETF) + dev branches"
Replace-SubText $body "I-D in progress" "I-D in progress – stable (IETF) + dev branches"

# ---------------------------------------------------------------------------
# Slide 3: three free-floating text boxes
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

$s3.Shapes.Item(8).TextFrame.TextRange.Text = "Nice Models, People can use, Both Syntax"
$s3.Shapes.Item(9).TextFrame.TextRange.Text = "Contract not to break compatibility, validation syntax"
$s3.Shapes.Item(10).TextFrame.TextRange.Text = "Framework Syntax ++"
